# Refresh cryptos.xlsx price (D) / 1h volume change (E) columns
# with the latest scraped values (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.774.36'
$ws.Range('E2').Value = '  +5.51%  '
$ws.Range('D3').Value = '2.224.23'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.23'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.76'
$ws.Range('E7').Value = '  -3.14%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.90'
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +5.36%  '
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '2.555.58'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.67'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.74'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D18').Value = '2.227.94'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('D19').Value = '41.613.32'
$ws.Range('E19').Value = '  +5.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.74'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.02'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.55'
$ws.Range('E23').Value = '  +9.57%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  +1.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.31'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.54'
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('E28').Value = '  +3.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.68'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.92'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.63'
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('E34').Value = '  +5.47%  '
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0624'
$ws.Range('E36').Value = '  +1.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.64'
$ws.Range('E37').Value = '  -4.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.68'
$ws.Range('E38').Value = '  -4.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.36'
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000256'
$ws.Range('E40').Value = '  +34.41%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +5.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.80'
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0977'
$ws.Range('E45').Value = '  +6.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.22'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.60'
$ws.Range('E47').Value = '  -3.76%  '
$ws.Range('D48').Value = '1.468.41'
$ws.Range('E48').Value = '  -2.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.54'
$ws.Range('E49').Value = '  -6.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.07'
$ws.Range('E51').Value = '  -1.40%  '
